# Update cryptos list: refreshed prices / 1h volume %, and restore the
# WEMIXToken / Stacks / NEARProtocol row order (rows 48-50).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.035.96"
$ws.Range("E2").Value = "  +3.16%  "
$ws.Range("D3").Value = "3.719.04"
$ws.Range("E3").Value = "  +6.51%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'420.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "'132.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("D7").Value = "3.712.38"
$ws.Range("E7").Value = "  +6.59%  "
$ws.Range("D8").Value = "'0.646"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("D11").Value = "'0.185"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.18%  "
$ws.Range("D12").Value = "'0.0000410"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +57.47%  "
$ws.Range("D13").Value = "'43.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "'10.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.55%  "
$ws.Range("D15").Value = "4.296.20"
$ws.Range("E15").Value = "  +6.01%  "
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "'20.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "3.726.04"
$ws.Range("E18").Value = "  +6.91%  "
$ws.Range("D19").Value = "'13.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.23%  "
$ws.Range("E20").Value = "  +4.45%  "
$ws.Range("D21").Value = "67.040.96"
$ws.Range("E21").Value = "  +3.35%  "
$ws.Range("D22").Value = "'446.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").Value = "'16.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +25.73%  "
$ws.Range("D24").Value = "'89.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'3.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "'37.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.37%  "
$ws.Range("D27").Value = "'10.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("D28").Value = "'3.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  +3.62%  "
$ws.Range("D30").Value = "'12.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("E31").Value = "  +10.04%  "
$ws.Range("E32").Value = "  +3.03%  "
$ws.Range("D33").Value = "'7.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.04%  "
$ws.Range("D34").Value = "'0.166"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "'41.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.94%  "
$ws.Range("D36").Value = "'57.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("D39").Value = "0.0₃0742"
$ws.Range("E39").Value = "  +6.37%  "
$ws.Range("D40").Value = "'3.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +32.72%  "
$ws.Range("D41").Value = "'0.150"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.72%  "
$ws.Range("D42").Value = "'28.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +33.07%  "
$ws.Range("D43").Value = "'0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("E44").Value = "  +5.07%  "
$ws.Range("D45").Value = "'3.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +34.08%  "
$ws.Range("D46").Value = "'147.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("E47").Value = "  +6.00%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'4.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.54%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.00%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.00%  "
$ws.Range("D51").Value = "'0.311"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
